$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Common (sheet2) / RM_01 (sheet3): Test_Script_Id value 22302 -> 22746
# ---------------------------------------------------------------------
$wsCommon = $wb.Worksheets.Item("Common")
$wsCommon.Range("B2").Value = 22746

$wsRM01 = $wb.Worksheets.Item("RM_01")
$wsRM01.Range("B2").Value = 22746

# ---------------------------------------------------------------------
# Credentials (sheet1): turn the stray 2-row tail (A3 "RM User" header-
# only row / A4 "hr" row) into one full credential row for a
# Recruitment admin user, and drop the now-unused row.
# ---------------------------------------------------------------------
$wsCred = $wb.Worksheets.Item("Credentials")

# Remove the extra trailing row (old row 4, "hr") first so row 3 is the
# last row left.
$wsCred.Rows.Item(4).Delete()

# Rebuild row 3 as a full credentials entry.
$wsCred.Range("A3").Value = "Recruitment"
$wsCred.Range("B3").Value = "WHaque@SEMPRANRGU"
$wsCred.Range("C3").Value = "!!Mar1983"
$wsCred.Range("D3").Value = "Yes"
$wsCred.Range("E3").Value = "This is an Admin user for Recruitment"

# Link the username cell, then strip the auto-applied "Hyperlink" cell
# style back off so the cell keeps the default formatting.
$wsCred.Hyperlinks.Add($wsCred.Range("B3"), "mailto:WHaque@SEMPRANRGU")
$wsCred.Range("B3").ClearFormats()

# The second data row (row 2) loses its border/fill formatting and gets
# a slightly tighter row height.
$wsCred.Range("B2:D2").ClearFormats()
$wsCred.Rows.Item(2).RowHeight = 14.25

# Column width tweaks on the credentials sheet (columns now hold the
# longer "Recruitment" / "WHaque@SEMPRANRGU" values).
$wsCred.Columns.Item(1).ColumnWidth = 12.5
$wsCred.Columns.Item(2).ColumnWidth = 22.5

# ---------------------------------------------------------------------
# Selections per sheet (applied before the final sheet activation so the
# activation order matches the target: Credentials ends up selected).
# ---------------------------------------------------------------------
$wsCommon.Range("C17").Select()
$wsRM01.Range("B8").Select()
$wsCred.Range("C5").Select()
